$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $ws.Range("AM$row").Value = "DBLF"
}

for ($row = 34; $row -le 65; $row++) {
    $ws.Range("AM$row").Value = "Ballance-Aware-DBLF"
}
